$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 2454.7917
$ws.Range("I19").Value = 2120.7693
$ws.Range("J19").Value = 2849.5454
$ws.Range("K19").Value = 2120.7693
$ws.Range("L19").Value = 2849.5454
$ws.Range("M19").Value = -1945.7693
$ws.Range("N19").Value = -3199.5454
$ws.Range("H28").Value = 1241.238
$ws.Range("I28").Value = 703.2
$ws.Range("J28").Value = 2586.3333
$ws.Range("K28").Value = 703.2
$ws.Range("L28").Value = 2586.3333
$ws.Range("M28").Value = -218.2
$ws.Range("N28").Value = -3556.3333
$ws.Range("H54").Value = 22222
$ws.Range("I54").Value = 22222
$ws.Range("K54").Value = 22222
$ws.Range("M54").Value = -21736
$ws.Range("H76").Value = 8332.916999999999
$ws.Range("I76").Value = 8719
$ws.Range("K76").Value = 8719
$ws.Range("M76").Value = -8404
$ws.Range("H79").Value = 8332.916999999999
$ws.Range("I79").Value = 8719
$ws.Range("K79").Value = 8719
$ws.Range("M79").Value = -7627
$ws.Range("H100").Value = 3599.8333
$ws.Range("J100").Value = 10000
$ws.Range("L100").Value = 10000
$ws.Range("N100").Value = -11082
$ws.Range("H111").Value = 5023.1333
$ws.Range("I111").Value = 3486.182
$ws.Range("J111").Value = 9249.75
$ws.Range("K111").Value = 10458.546
$ws.Range("L111").Value = 27749.25
$ws.Range("M111").Value = -7391.545999999998
$ws.Range("N111").Value = -33883.25
$ws.Range("H132").Value = 4465.2104
$ws.Range("I132").Value = 4581.0557
$ws.Range("K132").Value = 13743.1671
$ws.Range("M132").Value = -11213.1671
$ws.Range("H138").Value = 3696.7778
$ws.Range("I138").Value = 1461.625
$ws.Range("K138").Value = 4384.875
$ws.Range("M138").Value = 755.125
$ws.Range("H141").Value = 6146.6
$ws.Range("I141").Value = 6059.625
$ws.Range("K141").Value = 18178.875
$ws.Range("M141").Value = -12998.875

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13346.6
$ws.Range("I32").Value = 1494.8948
$ws.Range("J32").Value = 33817.727
$ws.Range("K32").Value = 1494.8948
$ws.Range("L32").Value = 33817.727
$ws.Range("M32").Value = -1207.8948
$ws.Range("N32").Value = -34391.727

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 4444.8696
$ws.Range("I99").Value = 4363.9375
$ws.Range("J99").Value = 4629.857
$ws.Range("K99").Value = 4363.9375
$ws.Range("L99").Value = 4629.857
$ws.Range("M99").Value = -2865.9375
$ws.Range("N99").Value = -7625.857

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 27885.379
$ws.Range("I16").Value = 22724.883
$ws.Range("J16").Value = 35196.082
$ws.Range("K16").Value = 22724.883
$ws.Range("L16").Value = 35196.082
$ws.Range("M16").Value = -22437.883
$ws.Range("N16").Value = -35770.082
$ws.Range("H22").Value = 3998
$ws.Range("I22").Value = 3997.5
$ws.Range("K22").Value = 3997.5
$ws.Range("M22").Value = -3647.5
$ws.Range("H31").Value = 7022.069
$ws.Range("I31").Value = 7126.1665
$ws.Range("K31").Value = 7126.1665
$ws.Range("M31").Value = -6831.1665
$ws.Range("H34").Value = 7022.069
$ws.Range("I34").Value = 7126.1665
$ws.Range("K34").Value = 7126.1665
$ws.Range("M34").Value = -6924.1665
$ws.Range("H44").Value = 20000
$ws.Range("J44").Value = 20000
$ws.Range("L44").Value = 20000
$ws.Range("N44").Value = -20884
$ws.Range("H58").Value = 5275.0557
$ws.Range("I58").Value = 5561.923
$ws.Range("K58").Value = 5561.923
$ws.Range("M58").Value = -5358.923
$ws.Range("H94").Value = 760
$ws.Range("I94").Value = 598.75
$ws.Range("J94").Value = 921.25
$ws.Range("K94").Value = 598.75
$ws.Range("L94").Value = 921.25
$ws.Range("M94").Value = -147.75
$ws.Range("N94").Value = -1823.25
$ws.Range("H113").Value = 27885.379
$ws.Range("I113").Value = 22724.883
$ws.Range("J113").Value = 35196.082
$ws.Range("K113").Value = 22724.883
$ws.Range("L113").Value = 35196.082
$ws.Range("M113").Value = -20554.883
$ws.Range("N113").Value = -39536.082
$ws.Range("H122").Value = 2749.25
$ws.Range("I122").Value = 2736.5386
$ws.Range("J122").Value = 2804.3333
$ws.Range("K122").Value = 8209.6158
$ws.Range("L122").Value = 8412.999899999999
$ws.Range("M122").Value = -5759.6158
$ws.Range("N122").Value = -13312.9999
$ws.Range("H134").Value = 3919.92
$ws.Range("I134").Value = 3335.0952
$ws.Range("K134").Value = 10005.2856
$ws.Range("M134").Value = -7470.285600000001
$ws.Range("H136").Value = 5275.0557
$ws.Range("I136").Value = 5561.923
$ws.Range("K136").Value = 16685.769
$ws.Range("M136").Value = -14135.769

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 97619310
$ws.Range("J4").Value = 8333722.5
$ws.Range("L4").Value = 25001167.5
$ws.Range("N4").Value = -25001391.5
$ws.Range("H34").Value = 5106.75
$ws.Range("J34").Value = 5525.5454
$ws.Range("L34").Value = 16576.6362
$ws.Range("N34").Value = -16744.6362
$ws.Range("H38").Value = 72.59999999999999
$ws.Range("J38").Value = 86
$ws.Range("L38").Value = 258
$ws.Range("N38").Value = -952
$ws.Range("H39").Value = 7673.8
$ws.Range("I39").Value = 2022.5
$ws.Range("J39").Value = 9086.625
$ws.Range("K39").Value = 6067.5
$ws.Range("L39").Value = 27259.875
$ws.Range("M39").Value = -5773.5
$ws.Range("N39").Value = -27847.875
$ws.Range("H55").Value = 8276.869000000001
$ws.Range("J55").Value = 8276.869000000001
$ws.Range("L55").Value = 24830.607
$ws.Range("N55").Value = -25184.607
$ws.Range("H131").Value = 1409.4166
$ws.Range("I131").Value = 1089.5
$ws.Range("K131").Value = 3268.5
$ws.Range("M131").Value = 1771.5
$ws.Range("H134").Value = 4475.6
$ws.Range("I134").Value = 597
$ws.Range("K134").Value = 1791
$ws.Range("M134").Value = 3279
$ws.Range("H140").Value = 2575.8
$ws.Range("I140").Value = 2421.48
$ws.Range("J140").Value = 3347.4
$ws.Range("K140").Value = 7264.440000000001
$ws.Range("L140").Value = 10042.2
$ws.Range("M140").Value = -2084.440000000001
$ws.Range("N140").Value = -20402.2

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 6201.8887
$ws.Range("I43").Value = 6201.8887
$ws.Range("K43").Value = 6201.8887
$ws.Range("M43").Value = -6050.8887
$ws.Range("H70").Value = 5777.2144
$ws.Range("J70").Value = 5799.75
$ws.Range("L70").Value = 5799.75
$ws.Range("N70").Value = -6339.75
$ws.Range("H73").Value = 5777.2144
$ws.Range("J73").Value = 5799.75
$ws.Range("L73").Value = 5799.75
$ws.Range("N73").Value = -7671.75
$ws.Range("H80").Value = 6840.7144
$ws.Range("J80").Value = 6698.75
$ws.Range("L80").Value = 6698.75
$ws.Range("N80").Value = -8694.75
$ws.Range("H83").Value = 6840.7144
$ws.Range("J83").Value = 6698.75
$ws.Range("L83").Value = 33493.75
$ws.Range("N83").Value = -43477.75
$ws.Range("H122").Value = 3138.6943
$ws.Range("I122").Value = 2103.6667
$ws.Range("J122").Value = 5208.75
$ws.Range("K122").Value = 6311.000100000001
$ws.Range("L122").Value = 15626.25
$ws.Range("M122").Value = -3861.000100000001
$ws.Range("N122").Value = -20526.25

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 8992.450999999999
$ws.Range("I40").Value = 6844.95
$ws.Range("J40").Value = 12897
$ws.Range("K40").Value = 6844.95
$ws.Range("L40").Value = 12897
$ws.Range("M40").Value = -6708.95
$ws.Range("N40").Value = -13169

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 4618.92
$ws.Range("I122").Value = 3793.2222
$ws.Range("J122").Value = 6742.143
$ws.Range("K122").Value = 11379.6666
$ws.Range("L122").Value = 20226.429
$ws.Range("M122").Value = -8929.6666
$ws.Range("N122").Value = -25126.429
